$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: Xbee part now sourced as XBP9B-XCUT-001 ---
# Set the hyperlink/url text in F5 first (reclaims the old shared string slot,
# which also causes the now-orphaned "602-1964-ND" string to be dropped and
# everything above it to re-index down by one -- matching the target diff).
$ws.Range("F5").Value = "https://www.digikey.ca/product-detail/en/digi-international/XBP9B-XCUT-001/602-1295-ND/3043287"

# Digikey part number in B5 (new shared string, appended after the URL).
$ws.Range("B5").Value = "XBP9B-XCUT-001"

# Match formatting used elsewhere in the "Digikey Part #" column (e.g. B21).
$ws.Range("B21").Copy()
$ws.Range("B5").PasteSpecial(-4122)

# Updated unit price -> updates the dependent formula in G5 automatically.
$ws.Range("E5").Value = 54.3

# --- Hyperlinks: drop the stale F5 hyperlink, keep F11/F19/F23 pointing at
#     the same targets they always did ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.digikey.ca/product-detail/en/torex-semiconductor-ltd/XC6201P302MR-G/893-1189-1-ND/3906799")
$ws.Hyperlinks.Add($ws.Range("F19"), "https://www.digikey.ca/products/en?keywords=602-1964-ND")
$ws.Hyperlinks.Add($ws.Range("F23"), "https://www.adafruit.com/product/1781")

# Re-adding hyperlinks reformats the target cells with a fresh "Hyperlink"
# style; restore the original look (font/underline/colour) by pulling the
# format back off F5, which keeps its own pre-existing hyperlink-style
# formatting even though the live hyperlink on it is gone.
$ws.Range("F5").Copy()
$ws.Range("F11").PasteSpecial(-4122)
$ws.Range("F23").PasteSpecial(-4122)

# --- Cursor position ---
$ws.Range("B23").Select()

$wb.Save()
